$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value() -eq "T/R1") {
        $cell.Value = "T"
    }
}
